# Update the "Last Updated" timestamp on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 12:09 PM"

# Insert a new leading row of stock data on the "Stock List" sheet, which
# pushes every existing data row down by one and drops the last (now
# out-of-range) row so the sheet keeps its original dimensions (A1:H76).
$stocks = $wb.Worksheets.Item("Stock List")
$stocks.Rows.Item(2).Insert()
$stocks.Rows.Item(77).Delete()

# The insert copies formatting down from the header row; clear it so the
# new row matches the plain (unstyled) look of the other data rows.
$stocks.Range("A2:H2").ClearFormats()

$stocks.Range("A2").Value = "📋"
$stocks.Range("B2").Value = "CAPTRU-RE1"
$stocks.Range("C2").Value = "CAPTRU-RE1"
$stocks.Range("D2").Value = 5.67
$stocks.Range("E2").Value = -11.9565
$stocks.Range("F2").Value = "N/A"
$stocks.Range("G2").Value = "N/A"
$stocks.Range("H2").Value = 0
